# The DefaultUSERS test-data sheet tracks an "isLoggedIn" boolean per user
# in column K. Flip the logged-in state that was swapped between the
# "test" scratch user (row 2) and the "max_mustermann" user (row 3):
#   K2 (row 2 - "test" user)            TRUE  -> FALSE
#   K3 (row 3 - "max_mustermann" user)  FALSE -> TRUE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("K2").Value = $false
$ws.Range("K3").Value = $true
